$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Preserve the original (unstyled) cell style so forcing text storage
# via NumberFormat does not leave a stray style index behind.
$origStyle = $ws.Range("D2").Style

# Force the Price/Volume columns to be stored as text so that values
# which look numeric (e.g. "1.000", "0.9998") are not silently
# converted into numbers by Excel.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "30.629.51"
$ws.Range("E2").Value = "  +0.55%  "
$ws.Range("D3").Value = "1.883.26"
$ws.Range("E3").Value = "  +0.26%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "249.62"
$ws.Range("E5").Value = "  +1.19%  "
$ws.Range("E6").Value = "  +0.03%  "
$ws.Range("D7").Value = "0.4757"
$ws.Range("E7").Value = "  -0.14%  "
$ws.Range("D8").Value = "0.2940"
$ws.Range("E8").Value = "  +1.30%  "
$ws.Range("D9").Value = "0.06541"
$ws.Range("E9").Value = "  +0.24%  "
$ws.Range("D10").Value = "22.04"
$ws.Range("E10").Value = "  +0.59%  "
$ws.Range("D11").Value = "0.07740"
$ws.Range("E11").Value = "  +0.12%  "
$ws.Range("D12").Value = "96.87"
$ws.Range("E12").Value = "  -0.44%  "
$ws.Range("D13").Value = "0.7387"
$ws.Range("E13").Value = "  +0.08%  "
$ws.Range("D14").Value = "1.881.74"
$ws.Range("E14").Value = "  +0.19%  "
$ws.Range("D15").Value = "5.248"
$ws.Range("E15").Value = "  +2.26%  "
$ws.Range("D16").Value = "274.94"
$ws.Range("E16").Value = "  +0.74%  "
$ws.Range("D17").Value = "30.722.47"
$ws.Range("E17").Value = "  +0.91%  "
$ws.Range("D18").Value = "13.18"
$ws.Range("E18").Value = "  -3.20%  "
$ws.Range("D19").Value = "0.000007554"
$ws.Range("E19").Value = "  -0.49%  "
$ws.Range("D20").Value = "0.9998"
$ws.Range("E20").Value = "  -0.03%  "
$ws.Range("D21").Value = "2.128.85"
$ws.Range("D22").Value = "5.352"
$ws.Range("E22").Value = "  +1.84%  "
$ws.Range("D23").Value = "1.001"
$ws.Range("E23").Value = "  +0.04%  "
$ws.Range("D24").Value = "6.237"
$ws.Range("E24").Value = "  +0.80%  "
$ws.Range("E25").Value = "  -1.02%  "
$ws.Range("D26").Value = "164.15"
$ws.Range("E26").Value = "  +0.05%  "
$ws.Range("D27").Value = "18.86"
$ws.Range("E27").Value = "  -0.08%  "
$ws.Range("D28").Value = "1.916"
$ws.Range("E28").Value = "  -1.28%  "
$ws.Range("D29").Value = "1.344"
$ws.Range("E29").Value = "  -1.96%  "
$ws.Range("D30").Value = "0.09730"
$ws.Range("E30").Value = "  -2.16%  "
$ws.Range("D31").Value = "1.505"
$ws.Range("D32").Value = "4.295"
$ws.Range("E32").Value = "  -0.44%  "
$ws.Range("D33").Value = "4.171"
$ws.Range("E33").Value = "  +2.51%  "
$ws.Range("D34").Value = "0.04890"
$ws.Range("E35").Value = "  +0.13%  "
$ws.Range("D36").Value = "0.6999"
$ws.Range("E36").Value = "  -0.09%  "
$ws.Range("D37").Value = "2.721"
$ws.Range("E37").Value = "  +0.24%  "
$ws.Range("D38").Value = "0.01916"
$ws.Range("E38").Value = "  +2.61%  "
$ws.Range("D39").Value = "2.789"
$ws.Range("D40").Value = "6.308"
$ws.Range("E40").Value = "  -0.52%  "
$ws.Range("D41").Value = "75.58"
$ws.Range("E41").Value = "  +6.52%  "
$ws.Range("D42").Value = "2.030"
$ws.Range("E42").Value = "  +4.21%  "
$ws.Range("D43").Value = "0.4250"
$ws.Range("E43").Value = "  +1.10%  "
$ws.Range("D44").Value = "0.8408"
$ws.Range("E44").Value = "  +0.44%  "
$ws.Range("D45").Value = "1.000"
$ws.Range("E45").Value = "  +0.03%  "
$ws.Range("D46").Value = "102.66"
$ws.Range("E46").Value = "  -0.14%  "
$ws.Range("D47").Value = "9.410"
$ws.Range("E47").Value = "  +1.76%  "
$ws.Range("D48").Value = "7.063"
$ws.Range("E48").Value = "  -0.39%  "
$ws.Range("D49").Value = "35.61"
$ws.Range("D50").Value = "916.32"
$ws.Range("E50").Value = "  -1.29%  "
$ws.Range("D51").Value = "0.05767"
$ws.Range("E51").Value = "  +2.14%  "

# Restore the original style on the whole column range.
$ws.Range("D2:E51").Style = $origStyle

